$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates driven by the crypto-price refresh diff.
# Numeric-looking text values (e.g. "548.66") are written with a temporary
# Text number format so Excel keeps them as text (matching the source feed,
# which stores every price/volume cell as inline text), then the style is
# reset back to Normal so no stray formatting is left behind.

$ws.Range('D2').Value = '61.854.15'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '2.459.77'
$ws.Range('E3').Value = '  -2.60%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.585'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.54%  '
$ws.Range('D9').Value = '2.459.17'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('E10').Value = '  -4.26%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('E13').Value = '  -4.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '2.903.22'
$ws.Range('E15').Value = '  -2.72%  '
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('D17').Value = '61.828.71'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '2.457.42'
$ws.Range('E18').Value = '  -2.55%  '
$ws.Range('E19').Value = '  -4.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '320.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('E24').Value = '  +5.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range('E26').Value = '  -8.24%  '
$ws.Range('D27').Value = '2.580.67'
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  -5.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '532.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.32%  '
$ws.Range('E33').Value = '  -5.30%  '
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.74'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.58%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  -5.10%  '
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '139.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.63%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '143.71'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.42%  '
$ws.Range('E47').Value = '  -2.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0532'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.11%  '
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0931'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.45%  '
